$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    # Force the value to be stored as text (matches source data where every
    # cell - even number/date-looking ones - is plain text), then drop the
    # "quote prefix" number format Excel applies for a leading apostrophe so
    # the cell keeps the workbook's default (unstyled) appearance.
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).Style = "Normal"
}

# --- Fill in previously-empty coordinate/operation/zone cells on row 76 ---
$ws.Range("M76").Value = -58.473179
$ws.Range("N76").Value = -34.629138
Set-TextCell "O76" "Devoto"
Set-TextCell "P76" "Capital Norte"

# --- Fill in previously-empty coordinate/operation/zone cells on row 77 ---
$ws.Range("M77").Value = -58.400188
$ws.Range("N77").Value = -34.583882
Set-TextCell "O77" "Recoleta"
Set-TextCell "P77" "Capital Sur"

# --- New row 78 ---
Set-TextCell "A78" "6377"
Set-TextCell "B78" "7/8/2025"
Set-TextCell "C78" "GUARDIA VIEJA 4377"
Set-TextCell "D78" "5"
Set-TextCell "E78" "808099347"
Set-TextCell "F78" "Optical Power"
Set-TextCell "G78" "Pendiente"
Set-TextCell "H78" "Picada"
Set-TextCell "I78" "1"
Set-TextCell "J78" "Cambio"
Set-TextCell "K78" "Sin equipos"
Set-TextCell "L78" "Pasante"
$ws.Range("M78").Value = -58.426322
$ws.Range("N78").Value = -34.600097
Set-TextCell "O78" "Almagro"
Set-TextCell "P78" "Capital Sur"

# --- New row 79 ---
Set-TextCell "A79" "6383"
Set-TextCell "B79" "7/8/2025"
Set-TextCell "C79" "FALCON, RAMON L.,CNEL. 1411"
Set-TextCell "D79" "6"
Set-TextCell "E79" "808099320"
Set-TextCell "F79" "Optical Power"
Set-TextCell "G79" "Pendiente"
Set-TextCell "H79" "Picada"
Set-TextCell "I79" "1"
Set-TextCell "J79" "Cambio"
Set-TextCell "K79" "Sin equipos"
Set-TextCell "L79" "Pasante"
$ws.Range("M79").Value = -58.448523
$ws.Range("N79").Value = -34.62452
Set-TextCell "O79" "Boedo"
Set-TextCell "P79" "Capital Sur"
